$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3829.1619
$ws.Range("I138").Value = 1839.3704
$ws.Range("J138").Value = 5139.512
$ws.Range("K138").Value = 5518.1112
$ws.Range("L138").Value = 15418.536
$ws.Range("M138").Value = -378.1112000000003
$ws.Range("N138").Value = -25698.536

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23995.092
$ws.Range("I32").Value = 3964.0469
$ws.Range("J32").Value = 130827.336
$ws.Range("K32").Value = 3964.0469
$ws.Range("L32").Value = 130827.336
$ws.Range("M32").Value = -3677.0469
$ws.Range("N32").Value = -131401.336
$ws.Range("H74").Value = 2246.3462
$ws.Range("I74").Value = 1130.4546
$ws.Range("J74").Value = 3064.6667
$ws.Range("K74").Value = 1130.4546
$ws.Range("L74").Value = 3064.6667
$ws.Range("M74").Value = -256.4546
$ws.Range("N74").Value = -4812.6667
$ws.Range("H77").Value = 2246.3462
$ws.Range("I77").Value = 1130.4546
$ws.Range("J77").Value = 3064.6667
$ws.Range("K77").Value = 5652.273
$ws.Range("L77").Value = 15323.3335
$ws.Range("M77").Value = -1284.273
$ws.Range("N77").Value = -24059.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 449
$ws.Range("I22").Value = 449
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 449
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -276
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1100.8334
$ws.Range("I16").Value = 745.5
$ws.Range("J16").Value = 1278.5
$ws.Range("K16").Value = 745.5
$ws.Range("L16").Value = 1278.5
$ws.Range("M16").Value = -458.5
$ws.Range("N16").Value = -1852.5
$ws.Range("H31").Value = 18364.918
$ws.Range("I31").Value = 48878.094
$ws.Range("J31").Value = 2345.5
$ws.Range("K31").Value = 48878.094
$ws.Range("L31").Value = 2345.5
$ws.Range("M31").Value = -48583.094
$ws.Range("N31").Value = -2935.5
$ws.Range("H34").Value = 18364.918
$ws.Range("I34").Value = 48878.094
$ws.Range("J34").Value = 2345.5
$ws.Range("K34").Value = 48878.094
$ws.Range("L34").Value = 2345.5
$ws.Range("M34").Value = -48676.094
$ws.Range("N34").Value = -2749.5
$ws.Range("H45").Value = 10000
$ws.Range("I45").Value = 12000
$ws.Range("J45").Value = 9666.667
$ws.Range("K45").Value = 12000
$ws.Range("L45").Value = 9666.667
$ws.Range("M45").Value = -11407
$ws.Range("N45").Value = -10852.667
$ws.Range("H113").Value = 1100.8334
$ws.Range("I113").Value = 745.5
$ws.Range("J113").Value = 1278.5
$ws.Range("K113").Value = 745.5
$ws.Range("L113").Value = 1278.5
$ws.Range("M113").Value = 1424.5
$ws.Range("N113").Value = -5618.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1983.2794
$ws.Range("I68").Value = 1347.92
$ws.Range("J68").Value = 2352.6743
$ws.Range("K68").Value = 4043.76
$ws.Range("L68").Value = 7058.0229
$ws.Range("M68").Value = -3232.76
$ws.Range("N68").Value = -8680.0229
$ws.Range("H71").Value = 1983.2794
$ws.Range("I71").Value = 1347.92
$ws.Range("J71").Value = 2352.6743
$ws.Range("K71").Value = 12131.28
$ws.Range("L71").Value = 21174.0687
$ws.Range("M71").Value = -8075.280000000001
$ws.Range("N71").Value = -29286.0687
$ws.Range("H88").Value = 12000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 12000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 36000
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -36856
$ws.Range("H91").Value = 12000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 12000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 36000
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -38964
$ws.Range("H131").Value = 1542.7609
$ws.Range("I131").Value = 1872.5
$ws.Range("J131").Value = 1511.3572
$ws.Range("K131").Value = 5617.5
$ws.Range("L131").Value = 4534.071599999999
$ws.Range("M131").Value = -577.5
$ws.Range("N131").Value = -14614.0716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 88820.375
$ws.Range("I70").Value = 138466.27
$ws.Range("J70").Value = 6077.222
$ws.Range("K70").Value = 138466.27
$ws.Range("L70").Value = 6077.222
$ws.Range("M70").Value = -138196.27
$ws.Range("N70").Value = -6617.222
$ws.Range("H73").Value = 88820.375
$ws.Range("I73").Value = 138466.27
$ws.Range("J73").Value = 6077.222
$ws.Range("K73").Value = 138466.27
$ws.Range("L73").Value = 6077.222
$ws.Range("M73").Value = -137530.27
$ws.Range("N73").Value = -7949.222
$ws.Range("H113").Value = 2999.125
$ws.Range("I113").Value = 6000
$ws.Range("J113").Value = 1998.8334
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 1998.8334
$ws.Range("M113").Value = -3830
$ws.Range("N113").Value = -6338.8334
$ws.Range("H132").Value = 4676.7144
$ws.Range("I132").Value = 4508
$ws.Range("J132").Value = 5295.3335
$ws.Range("K132").Value = 13524
$ws.Range("L132").Value = 15886.0005
$ws.Range("M132").Value = -10994
$ws.Range("N132").Value = -20946.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 640
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 640
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 640
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1230
$ws.Range("H27").Value = 640
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 640
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 640
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -854
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H34").Value = 7000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 7000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 7000
$ws.Range("N34").Value = -7344
$ws.Range("H40").Value = 252370
$ws.Range("I40").Value = 501000
$ws.Range("J40").Value = 3740
$ws.Range("K40").Value = 501000
$ws.Range("L40").Value = 3740
$ws.Range("M40").Value = -500864
$ws.Range("N40").Value = -4012
$ws.Range("H132").Value = 9120.8
$ws.Range("I132").Value = 11346.182
$ws.Range("J132").Value = 3001
$ws.Range("K132").Value = 34038.546
$ws.Range("L132").Value = 9003
$ws.Range("M132").Value = -31508.546
$ws.Range("N132").Value = -14063
$ws.Range("H136").Value = 3178.8572
$ws.Range("I136").Value = 2700.3635
$ws.Range("J136").Value = 4933.3335
$ws.Range("K136").Value = 8101.0905
$ws.Range("L136").Value = 14800.0005
$ws.Range("M136").Value = -5551.0905
$ws.Range("N136").Value = -19900.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2290713.8
$ws.Range("I62").Value = 8697432
$ws.Range("J62").Value = 2600
$ws.Range("K62").Value = 8697432
$ws.Range("L62").Value = 2600
$ws.Range("M62").Value = -8696808
$ws.Range("N62").Value = -3848
$ws.Range("H65").Value = 2290713.8
$ws.Range("I65").Value = 8697432
$ws.Range("J65").Value = 2600
$ws.Range("K65").Value = 43487160
$ws.Range("L65").Value = 13000
$ws.Range("M65").Value = -43484040
$ws.Range("N65").Value = -19240
$ws.Range("H132").Value = 9816.25
$ws.Range("I132").Value = 10504.429
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 31513.287
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -28983.287
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 1521.3846
$ws.Range("I136").Value = 998.5
$ws.Range("J136").Value = 2358
$ws.Range("K136").Value = 2995.5
$ws.Range("L136").Value = 7074
$ws.Range("M136").Value = -445.5
$ws.Range("N136").Value = -12174
